$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.238.33"
$ws.Range("E2").Value = "  +3.60%  "
$ws.Range("D3").Value = "1.590.44"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +8.69%  "
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("D12").Value = "1.818.16"
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("D13").Value = "1.583.48"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").Value = "28.303.18"
$ws.Range("E16").Value = "  +3.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.54%  "
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").Value = "1.402.95"
$ws.Range("E34").Value = "  -3.87%  "
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("E36").Value = "  -6.74%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("E39").Value = "  +8.09%  "
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("E43").Value = "  -3.56%  "
$ws.Range("E44").Value = "  +7.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.985"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").Value = "1.732.11"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.04%  "
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0525"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.14%  "
